$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B:F (voltage magnitudes), rows 2-25
$arrBF = New-Object 'object[,]' 24,5
$arrBF[0,0] = [double]1.02
$arrBF[0,1] = [double]1.053650188724247
$arrBF[0,2] = [double]1.050880821334305
$arrBF[0,3] = [double]1.059413228233382
$arrBF[0,4] = [double]1.068230967277935
$arrBF[1,0] = [double]1.02
$arrBF[1,1] = [double]1.055403096830601
$arrBF[1,2] = [double]1.052200904220148
$arrBF[1,3] = [double]1.061019939041431
$arrBF[1,4] = [double]1.070028090516584
$arrBF[2,0] = [double]1.02
$arrBF[2,1] = [double]1.056534558412087
$arrBF[2,2] = [double]1.053052377789826
$arrBF[2,3] = [double]1.062057226333644
$arrBF[2,4] = [double]1.071188772994389
$arrBF[3,0] = [double]1.02
$arrBF[3,1] = [double]1.057009572302049
$arrBF[3,2] = [double]1.053409698958797
$arrBF[3,3] = [double]1.062492749013974
$arrBF[3,4] = [double]1.071676216978874
$arrBF[4,0] = [double]1.02
$arrBF[4,1] = [double]1.057089291382631
$arrBF[4,2] = [double]1.053469657614496
$arrBF[4,3] = [double]1.06256584312968
$arrBF[4,4] = [double]1.071758031649946
$arrBF[5,0] = [double]1.02
$arrBF[5,1] = [double]1.056540908114905
$arrBF[5,2] = [double]1.053057154827065
$arrBF[5,3] = [double]1.062063047964242
$arrBF[5,4] = [double]1.071195288217289
$arrBF[6,0] = [double]1.02
$arrBF[6,1] = [double]1.054243178211906
$arrBF[6,2] = [double]1.051327516386071
$arrBF[6,3] = [double]1.05995672035259
$arrBF[6,4] = [double]1.068838772663169
$arrBF[7,0] = [double]1.02
$arrBF[7,1] = [double]1.050172275822986
$arrBF[7,2] = [double]1.048258476013771
$arrBF[7,3] = [double]1.056226444579034
$arrBF[7,4] = [double]1.064668986434209
$arrBF[8,0] = [double]1.02
$arrBF[8,1] = [double]1.047442613310521
$arrBF[8,2] = [double]1.046197558995674
$arrBF[8,3] = [double]1.053726253064014
$arrBF[8,4] = [double]1.061876601015297
$arrBF[9,0] = [double]1.02
$arrBF[9,1] = [double]1.046256700794697
$arrBF[9,2] = [double]1.045301484091683
$arrBF[9,3] = [double]1.052640299790608
$arrBF[9,4] = [double]1.060664295946991
$arrBF[10,0] = [double]1.02
$arrBF[10,1] = [double]1.045815588513748
$arrBF[10,2] = [double]1.044968075995465
$arrBF[10,3] = [double]1.052236409042712
$arrBF[10,4] = [double]1.06021349622734
$arrBF[11,0] = [double]1.02
$arrBF[11,1] = [double]1.04591023656188
$arrBF[11,2] = [double]1.045039618977315
$arrBF[11,3] = [double]1.052323068731135
$arrBF[11,4] = [double]1.060310217006289
$arrBF[12,0] = [double]1.02
$arrBF[12,1] = [double]1.046220250896827
$arrBF[12,2] = [double]1.04527393609288
$arrBF[12,3] = [double]1.052606924740112
$arrBF[12,4] = [double]1.0606270429133
$arrBF[13,0] = [double]1.02
$arrBF[13,1] = [double]1.046411179378482
$arrBF[13,2] = [double]1.045418231228064
$arrBF[13,3] = [double]1.052781748510856
$arrBF[13,4] = [double]1.060822183612774
$arrBF[14,0] = [double]1.02
$arrBF[14,1] = [double]1.047521233648695
$arrBF[14,2] = [double]1.046256949852759
$arrBF[14,3] = [double]1.053798252227901
$arrBF[14,4] = [double]1.061956989158399
$arrBF[15,0] = [double]1.02
$arrBF[15,1] = [double]1.048216471001872
$arrBF[15,2] = [double]1.046782060759585
$arrBF[15,3] = [double]1.054434969365424
$arrBF[15,4] = [double]1.062667958150532
$arrBF[16,0] = [double]1.02
$arrBF[16,1] = [double]1.048621611236734
$arrBF[16,2] = [double]1.047087994441821
$arrBF[16,3] = [double]1.05480603315569
$arrBF[16,4] = [double]1.063082348328344
$arrBF[17,0] = [double]1.02
$arrBF[17,1] = [double]1.04875968968329
$arrBF[17,2] = [double]1.047192250145544
$arrBF[17,3] = [double]1.054932502075423
$arrBF[17,4] = [double]1.063223593287069
$arrBF[18,0] = [double]1.02
$arrBF[18,1] = [double]1.048141918043778
$arrBF[18,2] = [double]1.046725758099458
$arrBF[18,3] = [double]1.0543666891163
$arrBF[18,4] = [double]1.062591709682929
$arrBF[19,0] = [double]1.02
$arrBF[19,1] = [double]1.04612897643023
$arrBF[19,2] = [double]1.045204951267176
$arrBF[19,3] = [double]1.052523350682128
$arrBF[19,4] = [double]1.06053375940247
$arrBF[20,0] = [double]1.02
$arrBF[20,1] = [double]1.044859812710288
$arrBF[20,2] = [double]1.044245478788726
$arrBF[20,3] = [double]1.051361358055322
$arrBF[20,4] = [double]1.059236968045566
$arrBF[21,0] = [double]1.02
$arrBF[21,1] = [double]1.045532962740023
$arrBF[21,2] = [double]1.044754428490688
$arrBF[21,3] = [double]1.051977643102263
$arrBF[21,4] = [double]1.059924700154893
$arrBF[22,0] = [double]1.02
$arrBF[22,1] = [double]1.048175606490227
$arrBF[22,2] = [double]1.046751199943749
$arrBF[22,3] = [double]1.054397543018279
$arrBF[22,4] = [double]1.062626164028713
$arrBF[23,0] = [double]1.02
$arrBF[23,1] = [double]1.051227404975
$arrBF[23,2] = [double]1.049054473723362
$arrBF[23,3] = [double]1.057193102788796
$arrBF[23,4] = [double]1.065749118204259
$ws.Range("B2:F25").Value = $arrBF

# Columns I:M (voltage magnitudes), rows 2-25
$arrIM = New-Object 'object[,]' 24,5
$arrIM[0,0] = [double]1.037365202626547
$arrIM[0,1] = [double]1.058666419243568
$arrIM[0,2] = [double]1.053633869327472
$arrIM[0,3] = [double]1.062142800196143
$arrIM[0,4] = [double]1.07093670449366
$arrIM[1,0] = [double]1.037738308848142
$arrIM[1,1] = [double]1.060066399894185
$arrIM[1,2] = [double]1.054765180772999
$arrIM[1,3] = [double]1.06356173162683
$arrIM[1,4] = [double]1.072547329010808
$arrIM[2,0] = [double]1.03797720051891
$arrIM[2,1] = [double]1.060969202151511
$arrIM[2,2] = [double]1.055493970113598
$arrIM[2,3] = [double]1.064477033107039
$arrIM[2,4] = [double]1.073586889217515
$arrIM[3,0] = [double]1.038077027817219
$arrIM[3,1] = [double]1.061348014806462
$arrIM[3,2] = [double]1.055799585778427
$arrIM[3,3] = [double]1.064861157119964
$arrIM[3,4] = [double]1.074023305725723
$arrIM[4,0] = [double]1.038093754018945
$arrIM[4,1] = [double]1.061411576916442
$arrIM[4,2] = [double]1.055850855266771
$arrIM[4,3] = [double]1.064925614354194
$arrIM[4,4] = [double]1.074096546312869
$arrIM[5,0] = [double]1.037978536779535
$arrIM[5,1] = [double]1.060974266699099
$arrIM[5,2] = [double]1.055498056767319
$arrIM[5,3] = [double]1.064482168402531
$arrIM[5,4] = [double]1.073592723033375
$arrIM[6,0] = [double]1.037491822301841
$arrIM[6,1] = [double]1.059140193805011
$arrIM[6,2] = [double]1.054016878895823
$arrIM[6,3] = [double]1.062622930227199
$arrIM[6,4] = [double]1.071481573259736
$arrIM[7,0] = [double]1.036614604654652
$arrIM[7,1] = [double]1.055884216794038
$arrIM[7,2] = [double]1.051381576857984
$arrIM[7,3] = [double]1.059324422739514
$arrIM[7,4] = [double]1.067740793578946
$arrIM[8,0] = [double]1.03601640630981
$arrIM[8,1] = [double]1.053696618266694
$arrIM[8,2] = [double]1.049607114222979
$arrIM[8,3] = [double]1.057109690804466
$arrIM[8,4] = [double]1.06523222365636
$arrIM[9,0] = [double]1.035754152347266
$arrIM[9,1] = [double]1.052745178605022
$arrIM[9,2] = [double]1.048834444573803
$arrIM[9,3] = [double]1.056146791689281
$arrIM[9,4] = [double]1.064142310185645
$arrIM[10,0] = [double]1.035656249671517
$arrIM[10,1] = [double]1.05239112690615
$arrIM[10,2] = [double]1.048546780585649
$arrIM[10,3] = [double]1.05578852695394
$arrIM[10,4] = [double]1.063736898008373
$arrIM[11,0] = [double]1.035677272370015
$arrIM[11,1] = [double]1.052467101601693
$arrIM[11,2] = [double]1.048608515546876
$arrIM[11,3] = [double]1.055865403390602
$arrIM[11,4] = [double]1.063823886339172
$arrIM[12,0] = [double]1.03574606970163
$arrIM[12,1] = [double]1.052715925815531
$arrIM[12,2] = [double]1.048810679712389
$arrIM[12,3] = [double]1.05611718976106
$arrIM[12,4] = [double]1.064108810405182
$arrIM[13,0] = [double]1.035788392956567
$arrIM[13,1] = [double]1.052869148770916
$arrIM[13,2] = [double]1.048935151907833
$arrIM[13,3] = [double]1.056272243565721
$arrIM[13,4] = [double]1.064284285542455
$arrIM[14,0] = [double]1.036033742796249
$arrIM[14,1] = [double]1.053759672574338
$arrIM[14,2] = [double]1.049658301841867
$arrIM[14,3] = [double]1.05717351173939
$arrIM[14,4] = [double]1.065304478600857
$arrIM[15,0] = [double]1.036186776250515
$arrIM[15,1] = [double]1.054317141387843
$arrIM[15,2] = [double]1.050110751056782
$arrIM[15,3] = [double]1.057737797732539
$arrIM[15,4] = [double]1.065943420640554
$arrIM[16,0] = [double]1.0362757267273
$arrIM[16,1] = [double]1.054641899732523
$arrIM[16,2] = [double]1.050374241432264
$arrIM[16,3] = [double]1.058066560352872
$arrIM[16,4] = [double]1.066315750152157
$arrIM[17,0] = [double]1.036306003869001
$arrIM[17,1] = [double]1.05475256587966
$arrIM[17,2] = [double]1.050464014700607
$arrIM[17,3] = [double]1.05817859659132
$arrIM[17,4] = [double]1.066442645233359
$arrIM[18,0] = [double]1.036170389446463
$arrIM[18,1] = [double]1.054257372128956
$arrIM[18,2] = [double]1.050062250609575
$arrIM[18,3] = [double]1.057677294165696
$arrIM[18,4] = [double]1.065874904996826
$arrIM[19,0] = [double]1.035725824162322
$arrIM[19,1] = [double]1.052642671198058
$arrIM[19,2] = [double]1.048751165720505
$arrIM[19,3] = [double]1.056043061660886
$arrIM[19,4] = [double]1.064024923227439
$arrIM[20,0] = [double]1.035443471731789
$arrIM[20,1] = [double]1.051623707668534
$arrIM[20,2] = [double]1.047923009711609
$arrIM[20,3] = [double]1.055012069630888
$arrIM[20,4] = [double]1.062858461608526
$arrIM[21,0] = [double]1.03559342251948
$arrIM[21,1] = [double]1.052164238581579
$arrIM[21,2] = [double]1.048362397381087
$arrIM[21,3] = [double]1.055558953206478
$arrIM[21,4] = [double]1.063477143609545
$arrIM[22,0] = [double]1.036177794899875
$arrIM[22,1] = [double]1.054284380531817
$arrIM[22,2] = [double]1.05008416715806
$arrIM[22,3] = [double]1.057704634284414
$arrIM[22,4] = [double]1.065905865362428
$arrIM[23,0] = [double]1.036843728702271
$arrIM[23,1] = [double]1.056728897247428
$arrIM[23,2] = [double]1.052065922405376
$arrIM[23,3] = [double]1.060179883316742
$arrIM[23,4] = [double]1.06871040782612
$ws.Range("I2:M25").Value = $arrIM
